# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from 2023-09-06 (45175) to 2023-09-14 (45183).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
